# Swap the "Workman" and "Shawmut Design and Construction" work-experience
# blocks, and swap the "Sloan Kettering" and "Sea Girt Recreation"
# extracurricular blocks. Every paragraph's style/run formatting stays the
# same; only the visible text changes, so this is implemented as a series
# of Find/Replace operations. Unique placeholder tokens are used so that a
# first replacement never gets clobbered by a later one searching for the
# text it just produced.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Work Experience: Workman <-> Shawmut Design and Construction ---
Replace-Text "Workman" "@@PLACEHOLDER1@@"
Replace-Text "dec - nov" "@@PLACEHOLDER2@@"
Replace-Text "Compose press releases highlighting attribtues of titles" "@@PLACEHOLDER3@@"
Replace-Text "Engage with media outlets to increase book publicity" "@@PLACEHOLDER4@@"
Replace-Text "Ensure authors remained in the loop on publicity efforts" "@@PLACEHOLDER5@@"

Replace-Text "Shawmut Design and Construction" "Workman"
Replace-Text "sept - oct" "dec - nov"
Replace-Text "Copyedit internal documents" "Compose press releases highlighting attribtues of titles"
Replace-Text "Copywrite external marketing materials sent to prospective clients" "Engage with media outlets to increase book publicity"
Replace-Text "Initiate rewrite of Business Development marketing template" "Ensure authors remained in the loop on publicity efforts"

Replace-Text "@@PLACEHOLDER1@@" "Shawmut Design and Construction"
Replace-Text "@@PLACEHOLDER2@@" "sept - oct"
Replace-Text "@@PLACEHOLDER3@@" "Copyedit internal documents"
Replace-Text "@@PLACEHOLDER4@@" "Copywrite external marketing materials sent to prospective clients"
Replace-Text "@@PLACEHOLDER5@@" "Initiate rewrite of Business Development marketing template"

# --- Extracurriculars: Sloan Kettering <-> Sea Girt Recreation ---
Replace-Text "Sloan Kettering" "@@PLACEHOLDER6@@"
Replace-Text "jul-aug" "@@PLACEHOLDER7@@"
Replace-Text "Ensure ER is clean and ready for use" "@@PLACEHOLDER8@@"
Replace-Text "Maintain patient comfort while administering care" "@@PLACEHOLDER9@@"
Replace-Text "Detail operation to accompanying family members" "@@PLACEHOLDER10@@"

Replace-Text "Sea Girt Recreation" "Sloan Kettering"
Replace-Text "jun - july" "jul-aug"
Replace-Text "Maintain order with children" "Ensure ER is clean and ready for use"
Replace-Text "Ensure sports equipment was ready-to-use each day" "Maintain patient comfort while administering care"
Replace-Text "Cooperate with needs of parents" "Detail operation to accompanying family members"

Replace-Text "@@PLACEHOLDER6@@" "Sea Girt Recreation"
Replace-Text "@@PLACEHOLDER7@@" "jun - july"
Replace-Text "@@PLACEHOLDER8@@" "Maintain order with children"
Replace-Text "@@PLACEHOLDER9@@" "Ensure sports equipment was ready-to-use each day"
Replace-Text "@@PLACEHOLDER10@@" "Cooperate with needs of parents"
